# Slide 42 ("Code Generation for CompoundStmt") - Grammar Rule bullet:
#   compoundStmt = "{" statement "}" .
# becomes
#   compoundStmt = "{" statements "}" .
#
# The run that used to hold `  = "{" statement "}" .` is split into three
# runs so that only the word "statement" -> "statements" portion is touched,
# leaving the " = " prefix and the `"}" .` suffix as their own runs.

$quote = [char]34

$p = $ppt.ActivePresentation

# Find the slide that contains the compoundStmt grammar rule (normally
# slide 42), rather than hard-coding the index, so the script is resilient
# to minor deck reshuffles.
$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $txt = $shape.TextFrame.TextRange.Text
            if ($txt -ne $null -and $txt.Contains('compoundStmt = ' + $quote + '{' + $quote + ' statement ' + $quote + '}' + $quote + ' .')) {
                $targetSlide = $slide
                $targetShape = $shape
            }
        }
    }
}

if ($targetShape -eq $null) {
    throw "Could not locate the shape containing the compoundStmt grammar rule"
}

$tr = $targetShape.TextFrame.TextRange

# Locate the exact substring we need to re-split/update within the shape's
# text, then convert the 0-based .NET string offset to the 1-based offset
# used by TextRange.Characters(start, length).
$oldChunk = ' = ' + $quote + '{' + $quote + ' statement ' + $quote + '}' + $quote + ' .'
$fullText = $tr.Text
$charIdx0 = $fullText.IndexOf($oldChunk)
if ($charIdx0 -lt 0) {
    throw "Could not locate the target grammar-rule text to edit"
}
$startPos = $charIdx0 + 1

# Run A: ' = '  (keep as-is, no text change needed)
# Run B: '"{" statement ' -> '"{" statements '
# Run C: '"}" .' (keep as-is, but gets its own run once Run B is edited)

$runBOldText = $quote + '{' + $quote + ' statement '
$runBNewText = $quote + '{' + $quote + ' statements '

$runBStart = $startPos + 3          # skip past ' = '
$runBLen = $runBOldText.Length

$runB = $tr.Characters($runBStart, $runBLen)
$runB.Text = $runBNewText
